$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.123.44"
$ws.Range("E2").Value = "  +5.83%  "
$ws.Range("D3").Value = "1.922.50"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.19"
$ws.Range("E5").Value = "  +3.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5180"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4054"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08483"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.129"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.85"
$ws.Range("E11").Value = "  +2.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.25"
$ws.Range("E12").Value = "  +9.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.368"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("D14").Value = "1.924.70"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.387"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.95"
$ws.Range("E17").Value = "  +5.16%  "
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06744"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  +3.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.064"
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").Value = "30.124.88"
$ws.Range("E23").Value = "  +5.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.31"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.200"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").Value = "2.147.76"
$ws.Range("E26").Value = "  +2.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.24"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.69"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.11"
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.082"
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.103"
$ws.Range("E33").Value = "  +5.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.661"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02525"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06629"
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2221"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.241"
$ws.Range("E38").Value = "  +4.60%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.029"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.214"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6587"
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.66"
$ws.Range("E43").Value = "  +5.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6192"
$ws.Range("E44").Value = "  +2.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.29"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.758"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.074"
$ws.Range("E47").Value = "  +3.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.246"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.90"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.71"
$ws.Range("E51").Value = "  +4.57%  "
